$d = $word.ActiveDocument

# Update the date/title line
$d.Content.Find.Execute("2023-09-14 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-15 Friday", 2) | Out-Null

# Update table cell values (positional, since some old values repeat)
$table = $d.Tables.Item(1)
$table.Cell(1,1).Range.Text = "61÷6=10, 1"  # was: 78÷6=13, 0
$table.Cell(1,2).Range.Text = "41÷3=13, 2"  # was: 41÷6=6, 5
$table.Cell(1,3).Range.Text = "11÷7=1, 4"  # was: 75÷2=37, 1
$table.Cell(1,4).Range.Text = "40÷8=5, 0"  # was: 64÷2=32, 0
$table.Cell(1,5).Range.Text = "39÷4=9, 3"  # was: 72÷3=24, 0
$table.Cell(5,1).Range.Text = "41÷8=5, 1"  # was: 72÷8=9, 0
$table.Cell(5,2).Range.Text = "21÷6=3, 3"  # was: 77÷5=15, 2
$table.Cell(5,3).Range.Text = "23÷7=3, 2"  # was: 90÷5=18, 0
$table.Cell(5,4).Range.Text = "41÷4=10, 1"  # was: 64÷5=12, 4
$table.Cell(5,5).Range.Text = "60÷9=6, 6"  # was: 81÷3=27, 0
$table.Cell(9,1).Range.Text = "24÷3=8, 0"  # was: 71÷5=14, 1
$table.Cell(9,2).Range.Text = "69÷5=13, 4"  # was: 48÷4=12, 0
$table.Cell(9,3).Range.Text = "97÷4=24, 1"  # was: 99÷4=24, 3
$table.Cell(9,4).Range.Text = "66÷2=33, 0"  # was: 30÷4=7, 2
$table.Cell(9,5).Range.Text = "11÷6=1, 5"  # was: 67÷3=22, 1
$table.Cell(13,1).Range.Text = "81÷6=13, 3"  # was: 96÷5=19, 1
$table.Cell(13,2).Range.Text = "52÷7=7, 3"  # was: 65÷4=16, 1
$table.Cell(13,3).Range.Text = "58÷6=9, 4"  # was: 99÷4=24, 3
$table.Cell(13,4).Range.Text = "41÷3=13, 2"  # was: 22÷2=11, 0
$table.Cell(13,5).Range.Text = "10÷5=2, 0"  # was: 59÷9=6, 5
$table.Cell(17,1).Range.Text = "79÷7=11, 2"  # was: 21÷7=3, 0
$table.Cell(17,2).Range.Text = "56÷8=7, 0"  # was: 35÷9=3, 8
$table.Cell(17,3).Range.Text = "92÷6=15, 2"  # was: 58÷3=19, 1
$table.Cell(17,4).Range.Text = "17÷4=4, 1"  # was: 97÷6=16, 1
$table.Cell(17,5).Range.Text = "87÷7=12, 3"  # was: 45÷6=7, 3
